$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.016.51"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3
$ws.Range("D3").Value = "2.737.04"
$ws.Range("E3").Value = "  -0.54%  "

# Row 5
$ws.Range("D5").Value = "569.19"
$ws.Range("E5").Value = "  -1.29%  "

# Row 6
$ws.Range("D6").Value = "159.50"
$ws.Range("E6").Value = "  +1.29%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  -1.59%  "

# Row 9
$ws.Range("E9").Value = "  -1.00%  "

# Row 10
$ws.Range("E10").Value = "  +4.61%  "

# Row 11
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  -1.44%  "

# Row 12
$ws.Range("E12").Value = "  -0.36%  "

# Row 13
$ws.Range("D13").Value = "3.220.30"

# Row 14
$ws.Range("D14").Value = "26.78"
$ws.Range("E14").Value = "  +0.16%  "

# Row 15
$ws.Range("D15").Value = "63.612.21"
$ws.Range("E15").Value = "  -0.42%  "

# Row 16
$ws.Range("E16").Value = "  -1.17%  "

# Row 17
$ws.Range("D17").Value = "2.741.63"
$ws.Range("E17").Value = "  -0.54%  "

# Row 18
$ws.Range("D18").Value = "12.09"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("E19").Value = "  -1.32%  "

# Row 20
$ws.Range("D20").Value = "353.96"
$ws.Range("E20").Value = "  -1.34%  "

# Row 21
$ws.Range("E21").Value = "  -2.57%  "

# Row 22
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("E23").Value = "  -5.24%  "

# Row 24
$ws.Range("D24").Value = "64.27"
$ws.Range("E24").Value = "  -2.77%  "

# Row 25
$ws.Range("E25").Value = "  +0.44%  "

# Row 26
$ws.Range("E26").Value = "  +0.18%  "

# Row 27
$ws.Range("D27").Value = "8.42"
$ws.Range("E27").Value = "  -0.31%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").Value = "  -1.51%  "

# Row 29
$ws.Range("E29").Value = "  +0.45%  "

# Row 30
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  +3.23%  "

# Row 31
$ws.Range("E31").Value = "  +8.34%  "

# Row 32
$ws.Range("D32").Value = "163.87"
$ws.Range("E32").Value = "  -2.96%  "

# Row 33
$ws.Range("E33").Value = "  -0.53%  "

# Row 34
$ws.Range("E34").Value = "  -1.39%  "

# Row 35
$ws.Range("E35").Value = "  +2.01%  "

# Row 37
$ws.Range("E37").Value = "  +1.20%  "

# Row 38
$ws.Range("E38").Value = "  -0.51%  "

# Row 39
$ws.Range("D39").Value = "350.31"
$ws.Range("E39").Value = "  +6.31%  "

# Row 40
$ws.Range("E40").Value = "  +3.17%  "

# Row 41
$ws.Range("E41").Value = "  -1.10%  "

# Row 42
$ws.Range("D42").Value = "38.69"
$ws.Range("E42").Value = "  -1.87%  "

# Row 43
$ws.Range("D43").Value = "21.95"
$ws.Range("E43").Value = "  +1.45%  "

# Row 44
$ws.Range("D44").Value = "21.15"
$ws.Range("E44").Value = "  -2.45%  "

# Row 45
$ws.Range("D45").Value = "0.0583"
$ws.Range("E45").Value = "  -1.31%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "134.56"
$ws.Range("E46").Value = "  -0.89%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.624"
$ws.Range("E47").Value = "  -1.53%  "

# Row 48
$ws.Range("E48").Value = "  -0.95%  "

# Row 49
$ws.Range("E49").Value = "  -2.42%  "

# Row 50
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.15%  "

# Row 51
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  +0.02%  "
